$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so values like "1.007" are not
# auto-converted to numbers by Excel's smart-entry parsing.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.862.75"
$ws.Range("E2").Value = "  -2.04%  "

$ws.Range("D3").Value = "1.831.37"
$ws.Range("E3").Value = "  -1.78%  "

$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").Value = "310.35"
$ws.Range("E5").Value = "  -1.50%  "

$ws.Range("E6").Value = "  +0.16%  "

$ws.Range("D7").Value = "0.4606"
$ws.Range("E7").Value = "  -0.72%  "

$ws.Range("D8").Value = "0.3665"
$ws.Range("E8").Value = "  -1.38%  "

$ws.Range("D9").Value = "0.07159"
$ws.Range("E9").Value = "  -2.62%  "

$ws.Range("D10").Value = "0.8762"
$ws.Range("E10").Value = "  -1.18%  "

$ws.Range("D11").Value = "0.07894"
$ws.Range("E11").Value = "  -0.45%  "

$ws.Range("D12").Value = "19.54"
$ws.Range("E12").Value = "  -1.90%  "

$ws.Range("D13").Value = "1.830.11"
$ws.Range("E13").Value = "  -1.77%  "

$ws.Range("D14").Value = "5.325"
$ws.Range("E14").Value = "  -1.54%  "

$ws.Range("D15").Value = "6.368"
$ws.Range("E15").Value = "  -3.57%  "

$ws.Range("D16").Value = "86.73"
$ws.Range("E16").Value = "  -6.08%  "

$ws.Range("D17").Value = "1.008"
$ws.Range("E17").Value = "  +0.24%  "

$ws.Range("D18").Value = "0.000008709"
$ws.Range("E18").Value = "  -2.06%  "

$ws.Range("E19").Value = "  +0.18%  "

$ws.Range("D20").Value = "26.896.90"
$ws.Range("E20").Value = "  -2.04%  "

$ws.Range("D21").Value = "14.43"
$ws.Range("E21").Value = "  -2.96%  "

$ws.Range("D22").Value = "4.990"
$ws.Range("E22").Value = "  -3.10%  "

$ws.Range("D23").Value = "10.42"
$ws.Range("E23").Value = "  -1.36%  "

$ws.Range("D24").Value = "1.982"
$ws.Range("E24").Value = "  +4.28%  "

$ws.Range("D25").Value = "150.67"
$ws.Range("E25").Value = "  -1.62%  "

$ws.Range("D26").Value = "18.19"
$ws.Range("E26").Value = "  -1.67%  "

$ws.Range("D27").Value = "1.962"
$ws.Range("E27").Value = "  -5.43%  "

$ws.Range("D28").Value = "113.19"
$ws.Range("E28").Value = "  -2.97%  "

$ws.Range("D29").Value = "4.913"
$ws.Range("E29").Value = "  -4.38%  "

$ws.Range("D30").Value = "0.08817"
$ws.Range("E30").Value = "  -0.80%  "

$ws.Range("D31").Value = "3.131"
$ws.Range("E31").Value = "  +3.42%  "

$ws.Range("D32").Value = "0.7512"
$ws.Range("E32").Value = "  -0.46%  "

$ws.Range("D33").Value = "4.452"
$ws.Range("E33").Value = "  -0.81%  "

$ws.Range("D34").Value = "1.127"
$ws.Range("E34").Value = "  -3.21%  "

$ws.Range("D35").Value = "2.547"
$ws.Range("E35").Value = "  -4.08%  "

$ws.Range("D36").Value = "1.086"
$ws.Range("E36").Value = "  +0.52%  "

$ws.Range("D37").Value = "0.01928"
$ws.Range("E37").Value = "  -1.70%  "

$ws.Range("D38").Value = "2.930"
$ws.Range("E38").Value = "  -2.27%  "

$ws.Range("D39").Value = "0.05120"

$ws.Range("D40").Value = "6.891"
$ws.Range("E40").Value = "  -3.64%  "

$ws.Range("D41").Value = "0.4953"
$ws.Range("E41").Value = "  -4.48%  "

$ws.Range("D42").Value = "0.1591"
$ws.Range("E42").Value = "  -3.07%  "

$ws.Range("D43").Value = "8.316"
$ws.Range("E43").Value = "  -0.28%  "

$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").Value = "1.007"
$ws.Range("E44").Value = "  +0.19%  "

$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "0.4663"
$ws.Range("E45").Value = "  -3.96%  "

$ws.Range("E46").Value = "  -2.02%  "

$ws.Range("E47").Value = "  -1.65%  "

$ws.Range("D48").Value = "1.605"
$ws.Range("E48").Value = "  -2.86%  "

$ws.Range("D49").Value = "0.06098"
$ws.Range("E49").Value = "  -2.43%  "

$ws.Range("D50").Value = "64.38"
$ws.Range("E50").Value = "  -2.10%  "

$ws.Range("D51").Value = "36.22"
$ws.Range("E51").Value = "  -2.76%  "

# Restore default "Normal" style on column D so no stray number format
# remains applied to the cells (matches original formatting).
$ws.Range("D2:D51").Style = "Normal"
